# Updated cryptos list on Tue Nov 14 06:17:15 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.673.43'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.62%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.063.89'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +1.26%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.17'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.32%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.668'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.94%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.87'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -5.64%  '

# Row 9
$ws.Range("E9").Value = '  -0.35%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.365'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.56%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0754'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.88%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.107'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.95%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.932'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +6.14%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.84'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.05%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.362.86'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +1.14%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.51'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.61%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.100.33'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.08%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.599.76'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.73%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.29'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -4.43%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.21'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.73%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0867'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.92%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '238.72'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.67%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.29'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.03%  '

# Row 24
$ws.Range("E24").Value = '  +0.06%  '

# Row 25
$ws.Range("E25").Value = '  -3.28%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.41'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.43%  '

# Row 27
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.14'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.70%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '164.94'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -1.91%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.28'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.02%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.123'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -0.76%  '

# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.13'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -5.68%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.20'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +8.53%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.51'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.77%  '

# Row 34
$ws.Range("E34").Value = '  -1.72%  '

# Row 35
$ws.Range("E35").Value = '  -0.02%  '

# Row 36
$ws.Range("E36").Value = '  -0.35%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.23'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.25%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0828'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.23%  '

# Row 39
$ws.Range("E39").Value = '  -2.91%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.87'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -5.55%  '

# Row 41
$ws.Range("E41").Value = '  -2.09%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.89'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -7.27%  '

# Row 43
$ws.Range("E43").Value = '  -2.11%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '94.74'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.19%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0911'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -4.74%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.410.13'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +9.33%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.07'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.23%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.59'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +13.47%  '

# Row 49
$ws.Range("E49").Value = '  +0.50%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.29'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.69%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.252.27'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.33%  '
